$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy cell formats (styles) from existing columns to new D,E,F,G columns ---
# Column A header/data already correct style (s=1 header, s=2 data) - untouched
# Column D (German short label) should mirror Column B style (s=1 header row, s=3 data rows)
# Columns E,F,G (definitions) should mirror Column C style (s=1 header row, s=4 data rows)
$ws.Range("B1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("C1").Copy()
$ws.Range("E1:G1").PasteSpecial(-4122)
$ws.Range("B2:B7").Copy()
$ws.Range("D2:D7").PasteSpecial(-4122)
$ws.Range("C2:C7").Copy()
$ws.Range("E2:G7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Write new-string cells in the precise order that reproduces the original authoring order ---
$ws.Range('C2').Value = 'Not Online'
$ws.Range('D2').Value = 'Nicht online'
$ws.Range('G2').Value = 'Offline bereitgestellter Dienst, für den es nicht einmal eine informative Webseite gibt.'
$ws.Range('C3').Value = 'Information'
$ws.Range('D3').Value = 'Informationen'
$ws.Range('G3').Value = 'Der Benutzer erhält Informationen über den administrativen Ablauf (z. B. Zweck, Bedingungen der Anfrage usw.) und eventuell über die Art und Weise, wie er zu erledigen ist (z. B. Stellen, Schalterzeiten).'
$ws.Range('C4').Value = 'One way interaction'
$ws.Range('D4').Value = 'Einweg-Interaktion'
$ws.Range('G4').Value = 'Zusätzlich zu den Informationen werden dem Benutzer die Formulare zur Anforderung des gewünschten Verwaltungsaktes / des Verfahrens zur Verfügung gestellt, die dann über herkömmliche Kanäle weitergeleitet werden muss (z. B. Formular zur Änderung der Residenz oder Formulare zur Selbstbescheinigung).'
$ws.Range('C5').Value = 'Bidirectional interaction'
$ws.Range('D5').Value = 'Bidirektionale Interaktion'
$ws.Range('G5').Value = 'Der Benutzer kann den gewünschten Verwaltungsakt / Vorgang initiieren (z. B. das Formular kann online ausgefüllt und gesendet werden) und es werden online nur die Übernahme der vom Benutzer eingegebenen Daten garantiert, nicht aber deren kontextuelle Verarbeitung.'
$ws.Range('C6').Value = 'Transaction'
$ws.Range('D6').Value = 'Transaktion'
$ws.Range('G6').Value = 'Der Benutzer kann den gewünschten Verwaltungsakt / Vorgang initiieren, indem er die erforderlichen Daten zur Verfügung stellt und die entsprechende Transaktion vollständig online durchführt, einschließlich der Zahlung der erwarteten Kosten.'
$ws.Range('C7').Value = 'Customization'
$ws.Range('D7').Value = 'Personalisierung'
$ws.Range('G7').Value = 'Zusätzlich zum gesamten Zyklus des Verwaltungsverfahrens von Interesse, den er selbst online ausführt, erhält der Benutzer Informationen (Erinnerung an die Fristen, Zurückgabe des Verfahrensergebnisses, etc.), die ihm auf der Besis des verknüpften Profils im Voraus gesendet werden (Pro-Aktivität)'
$ws.Range('A1').Value = 'codice_ 1_livello'
$ws.Range('B1').Value = 'label_ITA_1_Livello'
$ws.Range('C1').Value = 'label_ENG_1_livello'
$ws.Range('D1').Value = 'label_DEU_1_livello'
$ws.Range('F2').Value = 'Offline service for which it does not exist any informational web page'
$ws.Range('F3').Value = 'Users are informed about the administrative process that regards the service (e.g., objectives, how to require it, etc) and about the way in which the service can be used (e.g., opening hours of the information desk, location, etc)'
$ws.Range('F4').Value = 'In addition to the information, users can use online forms in order to start a request for an administrative act of interest. Afterwords, the user must send the forms via traditional channels'
$ws.Range('F5').Value = 'The user can start an administrative act of interest online (e.g., the form can be filled in and sent online) and it is guaranteed that the data is provided online, only; it is not guaranteed the concurrent online data processing.'
$ws.Range('F6').Value = 'The user can start an administrative act of interest online by providing the necessary data. The user can then carry out the transaction entirely online, including the possible payment of costs related to the service.'
$ws.Range('F7').Value = 'The user can carried out an entire administrative act online and (s)he can receive information according to his/her profile (e.g., it can receive information about deadlines, about the result of an administrative act, etc.) This is also called pro-active service.'
$ws.Range('E1').Value = 'definizione_ITA'
$ws.Range('F1').Value = 'definizione_ENG'
$ws.Range('G1').Value = 'definizione_DEU'

# --- Write the remaining brand-new cells (E2:E7) which reuse pre-existing shared strings ---
$ws.Range('E2').Value = 'Servizio erogato offline, per il quale non esiste nemmeno una pagina web informativa.'
$ws.Range('E3').Value = 'Sono fornite all''utente informazioni sul procedimento amministrativo (es. finalità, termini di richiesta, ecc.) ed eventualmente sulle modalità di espletamento (es. sedi, orari di sportello).'
$ws.Range('E4').Value = 'Oltre alle informazioni, sono resi disponibili all''utente i moduli per la richiesta dell''atto/procedimento amministrativo di interesse che dovrà poi essere inoltrata attraverso canali tradizionali (es. modulo di variazione residenza o moduli di autocertificazione).'
$ws.Range('E5').Value = 'L''utente può avviare l''atto/procedimento amministrativo di interesse (es. il modulo può essere compilato e inviato on line) e viene garantita on line solo la presa in carico dei dati immessi dall''utente e non la loro contestuale elaborazione.'
$ws.Range('E6').Value = 'L''utente può avviare l''atto/procedimento amministrativo di interesse fornendo i dati necessari ed eseguire la transazione corrispondente interamente on line, incluso l''eventuale pagamento dei costi previsti.'
$ws.Range('E7').Value = 'L''utente, oltre ad eseguire on line l''intero ciclo del procedimento amministrativo di interesse riceve informazioni (sono ricordate le scadenze, è restituito l''esito del procedimento, ecc.), che gli sono inviate preventivamente, sulla base del profilo collegato (c.d. pro-attività).'

# --- Row heights ---
$ws.Rows.Item(3).RowHeight = 90
$ws.Rows.Item(4).RowHeight = 105
$ws.Rows.Item(5).RowHeight = 105
$ws.Rows.Item(6).RowHeight = 90
$ws.Rows.Item(7).RowHeight = 105

# --- Column widths ---
$ws.Columns.Item(4).ColumnWidth = 21.1
$ws.Range("E1:F1").EntireColumn.ColumnWidth = 39.0
$ws.Columns.Item(7).ColumnWidth = 41.1

# --- View: zoom + selection ---
$excel.ActiveWindow.Zoom = 150
$ws.Range("G2").Select()
